# Apply the "8f_cost_burden_disease" update (feedback-meeting data refresh):
#  - rename the general_practitioner column header to "gp"
#  - insert a new "measuring_devices" cost column between vision_aids and
#    transportation (this pushes transportation .. productivity_absent one
#    column to the right, K:R -> L:S)
#  - set the measuring_devices value for the vi_sev data row to 0
#  - leave the selection where the author left it after the edit

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename "general_practitioner" -> "gp" (column B header, row 1)
$ws.Range("B1").Value = "gp"

# Insert a new column at K; everything from K (transportation) onward
# shifts right by one, so the new blank column becomes K and the old K
# (transportation) becomes L, etc.
$ws.Columns("K").Insert()

# Populate the newly inserted "measuring_devices" column.
$ws.Range("K1").Value = "measuring_devices"
$ws.Range("K2").Value = 0

# The trailing header cells (informal_care_personal .. productivity_absent,
# now P1:S1) lose their explicit header style after the insert shift.
$ws.Range("P1:S1").ClearFormats()

# Match the author's final selection/scroll position after the edit.
$ws.Range("O2").Select()
